$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NutritionalData")
$ws.Range("A232").Value = "row1"
$ws.Range("B232").Value = 480
$ws.Range("C232").Value = 23
$ws.Range("D232").Value = 4
$ws.Range("E232").Value = 7
$ws.Range("F232").Value = 63
$ws.Range("G232").Value = 6
$ws.Range("H232").Value = 370
